$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.535.12'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.827.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.84'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5187'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3886'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08406'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +9.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.117'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.428'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.15'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.523'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.821.95'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001130'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.56'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06612'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.071'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.562.78'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.42'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.13'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.032.88'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.413'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.77'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1094'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.101'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.741'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07460'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2221'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02369'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.228'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.804'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.52'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6337'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.190'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.401'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.59'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.783'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5986'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.02'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.991'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.204'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.65'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.82%  '
